# Bilag nummer tilfojet til titlen (og sprog-mærkning ryddet op i forbindelse
# hermed, som Word gør ved almindelig redigering).
$d = $word.ActiveDocument

function Set-ParaXml($paraIndex, $bodyInner) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`r`n" +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $bodyInner +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $p = $d.Paragraphs.Item($paraIndex)
    $p.Range.InsertXML($pkg)
}

# 1. Titel: "OC0302  - tjekEmne" -> "Bilag 30  - OC0302 tjekEmne"
$body1 = '<w:p><w:pPr><w:pStyle w:val="Titel"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Bilag 30  - </w:t></w:r>' +
    '<w:r><w:t>OC030</w:t></w:r>' +
    '<w:r><w:t>2</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>t</w:t></w:r>' +
    '<w:r><w:t>jekEmne</w:t></w:r>' +
    '</w:p>'
Set-ParaXml 1 $body1

# 2. "Operation:" overskrift faar engelsk sprogmarkering (en-US)
$body5 = '<w:p><w:pPr><w:pStyle w:val="Overskrift1"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>Operation:</w:t></w:r>' +
    '</w:p>'
Set-ParaXml 5 $body5

# 3. "tjekEmne(emne : String)" faar engelsk sprogmarkering (en-US)
$body6 = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>tjekEmne(</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>emne</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> : String</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
    '</w:p>'
Set-ParaXml 6 $body6

# 4. Tom linje derefter faar ogsaa engelsk sprogmarkering (en-US)
$body7 = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
Set-ParaXml 7 $body7

# 5. "UC03 Send Besked" samles i eet run
$body9 = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>UC03 Send Besked</w:t></w:r>' +
    '</w:p>'
Set-ParaXml 9 $body9

# 6. "Postconditions:" samles i eet run
$body15 = '<w:p><w:pPr><w:pStyle w:val="Overskrift1"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Postconditions:</w:t></w:r>' +
    '</w:p>'
Set-ParaXml 15 $body15

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output ($i.ToString() + ": [" + $p.Range.Text + "]")
}
